{"js": "const replacements = [\n  [\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\"],\n  [\"93\u00f73=\", \"40\u00f78=\"],\n  [\"68\u00f76=\", \"30\u00f73=\"],\n  [\"73\u00f73=\", \"89\u00f78=\"],\n  [\"67\u00f76=\", \"69\u00f75=\"],\n  [\"89\u00f75=\", \"85\u00f79=\"],\n  [\"27\u00f74=\", \"99\u00f77=\"],\n  [\"24\u00f78=\", \"49\u00f77=\"],\n  [\"33\u00f74=\", \"70\u00f72=\"],\n  [\"49\u00f79=\", \"61\u00f75=\"],\n  [\"64\u00f74=\", \"69\u00f74=\"],\n  [\"20\u00f73=\", \"67\u00f75=\"],\n  [\"73\u00f74=\", \"63\u00f73=\"],\n  [\"18\u00f76=\", \"83\u00f77=\"],\n  [\"61\u00f79=\", \"90\u00f72=\"],\n  [\"43\u00f72=\", \"52\u00f75=\"],\n  [\"31\u00f74=\", \"12\u00f77=\"],\n  [\"13\u00f76=\", \"39\u00f72=\"],\n  [\"38\u00f73=\", \"58\u00f73=\"],\n  [\"35\u00f79=\", \"59\u00f73=\"],\n  [\"68\u00f75=\", \"25\u00f76=\"],\n  [\"31\u00f77=\", \"55\u00f78=\"],\n  [\"74\u00f73=\", \"38\u00f75=\"],\n  [\"41\u00f76=\", \"43\u00f78=\"],\n  [\"35\u00f75=\", \"63\u00f75=\"],\n  [\"80\u00f77=\", \"81\u00f72=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-13 Wednesday\", \"2024-11-14 Thursday\"),\n    @(\"93\u00f73=\", \"40\u00f78=\"),\n    @(\"68\u00f76=\", \"30\u00f73=\"),\n    @(\"73\u00f73=\", \"89\u00f78=\"),\n    @(\"67\u00f76=\", \"69\u00f75=\"),\n    @(\"89\u00f75=\", \"85\u00f79=\"),\n    @(\"27\u00f74=\", \"99\u00f77=\"),\n    @(\"24\u00f78=\", \"49\u00f77=\"),\n    @(\"33\u00f74=\", \"70\u00f72=\"),\n    @(\"49\u00f79=\", \"61\u00f75=\"),\n    @(\"64\u00f74=\", \"69\u00f74=\"),\n    @(\"20\u00f73=\", \"67\u00f75=\"),\n    @(\"73\u00f74=\", \"63\u00f73=\"),\n    @(\"18\u00f76=\", \"83\u00f77=\"),\n    @(\"61\u00f79=\", \"90\u00f72=\"),\n    @(\"43\u00f72=\", \"52\u00f75=\"),\n    @(\"31\u00f74=\", \"12\u00f77=\"),\n    @(\"13\u00f76=\", \"39\u00f72=\"),\n    @(\"38\u00f73=\", \"58\u00f73=\"),\n    @(\"35\u00f79=\", \"59\u00f73=\"),\n    @(\"68\u00f75=\", \"25\u00f76=\"),\n    @(\"31\u00f77=\", \"55\u00f78=\"),\n    @(\"74\u00f73=\", \"38\u00f75=\"),\n    @(\"41\u00f76=\", \"43\u00f78=\"),\n    @(\"35\u00f75=\", \"63\u00f75=\"),\n    @(\"80\u00f77=\", \"81\u00f72=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $wdFindContinue = 1\n    $wdReplaceAll = 2\n    $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $newText,\n        $wdReplaceAll\n    ) | Out-Null\n}"}
